$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 3 (pushes old "Lead / Auto Loan" row down to row 4) ---
$ws.Rows.Item(3).Insert()

# --- New row 3: Case / Disputed Debit Card Transaction / <debit dispute description> ---
$ws.Range("A3").Value = "Case"
$ws.Range("B3").Value = "Disputed Debit Card Transaction"
$ws.Range("C3").Value = "`n    <h3><u>Debit Card Dispute</u></h3>`n    1. Immediately block this card to prevent any further fraudulent charges.`n    <br>`n      <br>`n        2. Before submitting the dispute, ensure the customer/member has already called the merchant to inquire about the charge.`n        <br>`n          <br>`n            3. Inform customer/member it will take up to 10 business days to receive provisional credit for this charge`n            <br>`n              <br>`n                <strong>Links</strong>`n                <ul>`n                  <li><a href=`"https://crmnext.us`" target=`"_blank`"> Policy and Procedures Manual</a></li>`n                  <li><a href=`"https://crmnext.us`" target=`"_blank`"> Policy and Procedures Manual</a></li>`n                  <li><a href=`"https://crmnext.us`" target=`"_blank`"> Policy and Procedures Manual</a></li>`n                  <li><a href=`"https://crmnext.us`" target=`"_blank`"> Policy and Procedures Manual</a></li>`n                </ul>"
$ws.Range("C3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 404

# --- Row 4 (formerly row 3): update description cell, add wrap style, set height ---
$ws.Range("C4").Value = "<h3><u>Auto Loan</u></h3>`n<br>`n<br>`nSteps to submit an auto loan!"
$ws.Range("C4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 68

# --- Re-establish the sort memory (AutoFilter sort state) over the data range ---
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A1:A4"))
$sort.SetRange($ws.Range("A1:C4"))
$sort.Header = 1
$sort.Apply()

# --- Update the view: scroll so row 2 is at top, select C5 ---
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("C5").Select()
